# CUS15: actualizacion de servicios, scripts y archivos de cotizacion
# Adds the quotation number label/value to the COTIZACION sheet header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add "Cotizacion N°" label (C1) and its value (D1) ---
# Copy A3's formatting (bold, blue sub-header style used elsewhere on the
# sheet for section headers) onto C1 so the new label matches the rest of
# the form, then set the real text/value for the new cells.
$ws.Range("A3").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("D1").Value = 600231

# --- Restore the current selection state to the cell the user ended up on ---
$ws.Activate()
$ws.Range("D5").Select()
